$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.080.05"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.652.77"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("D5").Value = "218.25"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "0.5310"
$ws.Range("E6").Value = "  +1.76%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "0.2615"
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("D9").Value = "0.06295"
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("E10").Value = "  -3.80%  "
$ws.Range("D11").Value = "0.07739"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.683.88"
$ws.Range("E12").Value = "  +1.71%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.481"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "0.5455"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").Value = "0.0₅8110"
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("D16").Value = "65.23"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "26.110.47"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "4.561"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").Value = "193.89"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").Value = "10.04"
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("D22").Value = "5.997"
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").Value = "139.77"
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("D25").Value = "0.1244"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").Value = "7.259"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "16.20"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "1.432"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").Value = "0.05923"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").Value = "3.508"
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("D32").Value = "3.247"
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("D33").Value = "1.544"
$ws.Range("E33").Value = "  -6.47%  "
$ws.Range("D34").Value = "2.415"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").Value = "0.9431"
$ws.Range("E35").Value = "  -4.16%  "
$ws.Range("D36").Value = "2.755"
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("D37").Value = "0.5665"
$ws.Range("E37").Value = "  -4.60%  "
$ws.Range("D38").Value = "0.01608"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").Value = "5.849"
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("D40").Value = "0.8475"
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "1.009.43"
$ws.Range("E42").Value = "  -3.08%  "
$ws.Range("D43").Value = "100.65"
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("D44").Value = "1.795.18"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "57.05"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("D46").Value = "0.0₈106"
$ws.Range("E46").Value = "  -2.81%  "
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").Value = "0.4295"
$ws.Range("E48").Value = "  +1.63%  "
$ws.Range("D49").Value = "1.479"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05154"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.846"
$ws.Range("E51").Value = "  -3.68%  "
